# Auto-generated edit script applying the diff to before.xlsx
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 91736.836
$ws.Range("I11").Value = 91736.836
$ws.Range("K11").Value = 91736.836
$ws.Range("M11").Value = -91596.836
$ws.Range("H19").Value = 16787.334
$ws.Range("J19").Value = 28685.572
$ws.Range("L19").Value = 28685.572
$ws.Range("N19").Value = -29035.572
$ws.Range("H32").Value = 1601.3
$ws.Range("I32").Value = 885
$ws.Range("J32").Value = 1780.375
$ws.Range("K32").Value = 885
$ws.Range("L32").Value = 1780.375
$ws.Range("M32").Value = -559
$ws.Range("N32").Value = -2432.375
$ws.Range("H43").Value = 1507.1428
$ws.Range("I43").Value = 1280
$ws.Range("J43").Value = 2075
$ws.Range("K43").Value = 1280
$ws.Range("L43").Value = 2075
$ws.Range("M43").Value = -1211
$ws.Range("N43").Value = -2213
$ws.Range("H98").Value = 6602.933
$ws.Range("I98").Value = 6110.7085
$ws.Range("J98").Value = 8571.833000000001
$ws.Range("K98").Value = 6110.7085
$ws.Range("L98").Value = 8571.833000000001
$ws.Range("M98").Value = -4612.7085
$ws.Range("N98").Value = -11567.833
$ws.Range("H108").Value = 49000
$ws.Range("J108").Value = 49000
$ws.Range("L108").Value = 49000
$ws.Range("N108").Value = -56680
$ws.Range("H122").Value = 6602.933
$ws.Range("I122").Value = 6110.7085
$ws.Range("J122").Value = 8571.833000000001
$ws.Range("K122").Value = 18332.1255
$ws.Range("L122").Value = 25715.499
$ws.Range("M122").Value = -15882.1255
$ws.Range("N122").Value = -30615.499
$ws.Range("H132").Value = 20003314
$ws.Range("I132").Value = 21279302
$ws.Range("K132").Value = 63837906
$ws.Range("M132").Value = -63835376
$ws.Range("H137").Value = 3999.9333
$ws.Range("I137").Value = 3100.3157
$ws.Range("J137").Value = 5553.8184
$ws.Range("K137").Value = 9300.947100000001
$ws.Range("L137").Value = 16661.4552
$ws.Range("M137").Value = -6750.947100000001
$ws.Range("N137").Value = -21761.4552
$ws.Range("H138").Value = 2125.9517
$ws.Range("I138").Value = 1751.0312
$ws.Range("J138").Value = 2525.8667
$ws.Range("K138").Value = 5253.0936
$ws.Range("L138").Value = 7577.6001
$ws.Range("M138").Value = -113.0936000000002
$ws.Range("N138").Value = -17857.6001
$ws.Range("H141").Value = 3205.5715
$ws.Range("I141").Value = 3158.44
$ws.Range("J141").Value = 3598.3333
$ws.Range("K141").Value = 9475.32
$ws.Range("L141").Value = 10794.9999
$ws.Range("M141").Value = -4295.32
$ws.Range("N141").Value = -21154.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 4701.5
$ws.Range("I10").Value = 4701.5
$ws.Range("K10").Value = 4701.5
$ws.Range("M10").Value = -4531.5
$ws.Range("H74").Value = 41669936
$ws.Range("I74").Value = 55557844
$ws.Range("J74").Value = 6207.3335
$ws.Range("K74").Value = 55557844
$ws.Range("L74").Value = 6207.3335
$ws.Range("M74").Value = -55556970
$ws.Range("N74").Value = -7955.3335
$ws.Range("H77").Value = 41669936
$ws.Range("I77").Value = 55557844
$ws.Range("J77").Value = 6207.3335
$ws.Range("K77").Value = 277789220
$ws.Range("L77").Value = 31036.6675
$ws.Range("M77").Value = -277784852
$ws.Range("N77").Value = -39772.6675
$ws.Range("H97").Value = 2501.5
$ws.Range("J97").Value = 3594
$ws.Range("L97").Value = 3594
$ws.Range("N97").Value = -4586
$ws.Range("H122").Value = 4440.1562
$ws.Range("I122").Value = 4021.7144
$ws.Range("J122").Value = 5239
$ws.Range("K122").Value = 12065.1432
$ws.Range("L122").Value = 15717
$ws.Range("M122").Value = -9615.143199999999
$ws.Range("N122").Value = -20617
$ws.Range("H139").Value = 220178.75
$ws.Range("J139").Value = 220178.75
$ws.Range("L139").Value = 220178.75
$ws.Range("N139").Value = -230458.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 153
$ws.Range("I22").Value = 66.25
$ws.Range("K22").Value = 66.25
$ws.Range("M22").Value = 106.75
$ws.Range("H86").Value = 2733.077
$ws.Range("I86").Value = 2150.842
$ws.Range("J86").Value = 4313.4287
$ws.Range("K86").Value = 2150.842
$ws.Range("L86").Value = 4313.4287
$ws.Range("M86").Value = -1027.842
$ws.Range("N86").Value = -6559.4287
$ws.Range("H89").Value = 2733.077
$ws.Range("I89").Value = 2150.842
$ws.Range("J89").Value = 4313.4287
$ws.Range("K89").Value = 10754.21
$ws.Range("L89").Value = 21567.1435
$ws.Range("M89").Value = -5138.210000000001
$ws.Range("N89").Value = -32799.14350000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3127.6897
$ws.Range("I31").Value = 2462.45
$ws.Range("J31").Value = 4606
$ws.Range("K31").Value = 2462.45
$ws.Range("L31").Value = 4606
$ws.Range("M31").Value = -2167.45
$ws.Range("N31").Value = -5196
$ws.Range("H34").Value = 3127.6897
$ws.Range("I34").Value = 2462.45
$ws.Range("J34").Value = 4606
$ws.Range("K34").Value = 2462.45
$ws.Range("L34").Value = 4606
$ws.Range("M34").Value = -2260.45
$ws.Range("N34").Value = -5010
$ws.Range("H122").Value = 2111.8823
$ws.Range("I122").Value = 2140.1333
$ws.Range("K122").Value = 6420.3999
$ws.Range("M122").Value = -3970.3999
$ws.Range("H134").Value = 4701.7075
$ws.Range("I134").Value = 4198.7354
$ws.Range("K134").Value = 12596.2062
$ws.Range("M134").Value = -10061.2062

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 800
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 800
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 2400
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -2624
$ws.Range("H15").Value = 37.25
$ws.Range("I15").Value = 26.6
$ws.Range("J15").Value = 55
$ws.Range("K15").Value = 79.80000000000001
$ws.Range("L15").Value = 165
$ws.Range("M15").Value = 60.19999999999999
$ws.Range("N15").Value = -445
$ws.Range("H16").Value = 272.33334
$ws.Range("J16").Value = 356.2
$ws.Range("L16").Value = 1068.6
$ws.Range("N16").Value = -1414.6
$ws.Range("H39").Value = 5020.2354
$ws.Range("I39").Value = 25000
$ws.Range("J39").Value = 3771.5
$ws.Range("K39").Value = 75000
$ws.Range("L39").Value = 11314.5
$ws.Range("M39").Value = -74706
$ws.Range("N39").Value = -11902.5
$ws.Range("H55").Value = 3249.8076
$ws.Range("J55").Value = 3625.8696
$ws.Range("L55").Value = 10877.6088
$ws.Range("N55").Value = -11231.6088
$ws.Range("H132").Value = 2637.4417
$ws.Range("I132").Value = 1763.4166
$ws.Range("J132").Value = 2798.8
$ws.Range("K132").Value = 15870.7494
$ws.Range("L132").Value = 25189.2
$ws.Range("M132").Value = -13340.7494
$ws.Range("N132").Value = -30249.2
$ws.Range("H133").Value = 6670.625
$ws.Range("I133").Value = 4506.7856
$ws.Range("J133").Value = 9700
$ws.Range("K133").Value = 13520.3568
$ws.Range("L133").Value = 29100
$ws.Range("M133").Value = -8460.356800000001
$ws.Range("N133").Value = -39220
$ws.Range("H139").Value = 1760690.2
$ws.Range("I139").Value = 2090304
$ws.Range("K139").Value = 6270912
$ws.Range("M139").Value = -6265772

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 9998.5
$ws.Range("I19").Value = 9998.5
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 9998.5
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -9710.5
$ws.Range("N19").ClearContents()
$ws.Range("H70").Value = 9102.73
$ws.Range("I70").Value = 8038.5884
$ws.Range("K70").Value = 8038.5884
$ws.Range("M70").Value = -7768.5884
$ws.Range("H73").Value = 9102.73
$ws.Range("I73").Value = 8038.5884
$ws.Range("K73").Value = 8038.5884
$ws.Range("M73").Value = -7102.5884
$ws.Range("H132").Value = 4648.6665
$ws.Range("I132").Value = 4644.25
$ws.Range("J132").Value = 4657.5
$ws.Range("K132").Value = 13932.75
$ws.Range("L132").Value = 13972.5
$ws.Range("M132").Value = -11402.75
$ws.Range("N132").Value = -19032.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 13000
$ws.Range("I42").Value = 13000
$ws.Range("K42").Value = 13000
$ws.Range("M42").Value = -12437
$ws.Range("H49").Value = 13000
$ws.Range("I49").Value = 13000
$ws.Range("K49").Value = 13000
$ws.Range("M49").Value = -12853
$ws.Range("H55").Value = 930.5294
$ws.Range("I55").Value = 562.2222
$ws.Range("K55").Value = 562.2222
$ws.Range("M55").Value = -389.2222
$ws.Range("H68").Value = 4222.5
$ws.Range("I68").Value = 2682.8572
$ws.Range("J68").Value = 15000
$ws.Range("K68").Value = 2682.8572
$ws.Range("L68").Value = 15000
$ws.Range("M68").Value = -1933.8572
$ws.Range("N68").Value = -16498
$ws.Range("H71").Value = 4222.5
$ws.Range("I71").Value = 2682.8572
$ws.Range("J71").Value = 15000
$ws.Range("K71").Value = 13414.286
$ws.Range("L71").Value = 75000
$ws.Range("M71").Value = -9670.286
$ws.Range("N71").Value = -82488

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 20000
$ws.Range("I86").Value = 20000
$ws.Range("K86").Value = 20000
$ws.Range("M86").Value = -18877
$ws.Range("H89").Value = 20000
$ws.Range("I89").Value = 20000
$ws.Range("K89").Value = 100000
$ws.Range("M89").Value = -94384
$ws.Range("H92").Value = 79500
$ws.Range("J92").Value = 79500
$ws.Range("L92").Value = 79500
$ws.Range("N92").Value = -84492
$ws.Range("H132").Value = 23892.95
$ws.Range("I132").Value = 9254.3125
$ws.Range("J132").Value = 82447.5
$ws.Range("K132").Value = 27762.9375
$ws.Range("L132").Value = 247342.5
$ws.Range("M132").Value = -25232.9375
$ws.Range("N132").Value = -252402.5
